$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# The "Reduziert" price tier was removed - rows that were priced "Reduziert"
# (Hungernde-Mitarbeiter in F7, Admin in F10) now use "Normal " pricing instead.
$ws.Range("F7").Value = "Normal "
$ws.Range("F10").Value = "Normal "

# Reflect the view state change captured in the saved workbook (scroll/selection).
$excel.ActiveWindow.ScrollColumn = 8
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("K13").Select()
